# Reorganización completa: limpieza de módulos antiguos, nuevas entregas y optimización
#
# - Renombra la hoja "Datos" a "sector"
# - Quita el formato de cabecera (negrita blanca sobre relleno azul)
# - Renueva los encabezados de columna (minúsculas) e inserta una nueva
#   columna "finca" entre "nombre" y "descripcion", agregando "comentario"
#   al final.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "sector"

# Drop the old bold/white-on-blue header styling from the existing header row
$ws.Range("A1:D1").ClearFormats()

# Rewrite the header row: codigo | nombre | finca | descripcion | comentario
# (a new "finca" column is inserted before the old "Descripcion" column, and
# "comentario" is appended as a brand-new trailing column)
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "finca"
$ws.Range("D1").Value = "descripcion"
$ws.Range("E1").Value = "comentario"
